$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.496.07"
$ws.Range("E2").Value = "  -2.64%  "
$ws.Range("D3").Value = "2.297.40"
$ws.Range("E3").Value = "  -3.58%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'538.58"
$ws.Range("E5").Value = "  -2.12%  "
$ws.Range("D6").Value = "'127.80"
$ws.Range("E6").Value = "  -5.09%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.566"
$ws.Range("E8").Value = "  -4.20%  "
$ws.Range("D9").Value = "2.295.67"
$ws.Range("E9").Value = "  -3.54%  "
$ws.Range("D10").Value = "'0.100"
$ws.Range("E10").Value = "  -2.04%  "
$ws.Range("D11").Value = "'5.47"
$ws.Range("E11").Value = "  -1.97%  "
$ws.Range("D12").Value = "'0.149"
$ws.Range("E12").Value = "  -1.09%  "
$ws.Range("D13").Value = "'0.329"
$ws.Range("E13").Value = "  -3.70%  "
$ws.Range("D14").Value = "2.706.33"
$ws.Range("E14").Value = "  -3.57%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").Value = "'22.96"
$ws.Range("E15").Value = "  -6.00%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "59.357.90"
$ws.Range("E16").Value = "  -2.63%  "
$ws.Range("E17").Value = "  -3.36%  "
$ws.Range("D18").Value = "2.290.26"
$ws.Range("E18").Value = "  -3.41%  "
$ws.Range("D19").Value = "'10.34"
$ws.Range("E19").Value = "  -4.64%  "
$ws.Range("D20").Value = "'4.00"
$ws.Range("E20").Value = "  -5.85%  "
$ws.Range("D21").Value = "'308.02"
$ws.Range("E21").Value = "  -3.73%  "
$ws.Range("D22").Value = "'6.46"
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  -0.39%  "
$ws.Range("D24").Value = "'62.76"
$ws.Range("E24").Value = "  -1.66%  "
$ws.Range("E25").Value = "  -3.63%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").Value = "'7.64"
$ws.Range("E27").Value = "  -6.84%  "
$ws.Range("D28").Value = "'1.33"
$ws.Range("E28").Value = "  -3.44%  "
$ws.Range("D29").Value = "'171.40"
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("E30").Value = "  +1.34%  "
$ws.Range("E31").Value = "  -3.12%  "
$ws.Range("D32").Value = "0.0₃0707"
$ws.Range("E32").Value = "  -6.61%  "
$ws.Range("D33").Value = "'5.74"
$ws.Range("E33").Value = "  -4.22%  "
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("E35").Value = "  -3.47%  "
$ws.Range("D36").Value = "'1.32"
$ws.Range("E36").Value = "  -7.71%  "
$ws.Range("D37").Value = "'17.62"
$ws.Range("E37").Value = "  -2.99%  "
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("D39").Value = "'3.95"
$ws.Range("E39").Value = "  -6.58%  "
$ws.Range("D40").Value = "'308.13"
$ws.Range("E40").Value = "  -5.97%  "
$ws.Range("D41").Value = "'37.58"
$ws.Range("E41").Value = "  -2.35%  "
$ws.Range("D42").Value = "'1.49"
$ws.Range("E42").Value = "  -6.44%  "
$ws.Range("D43").Value = "'134.90"
$ws.Range("E43").Value = "  -8.04%  "
$ws.Range("E44").Value = "  -3.78%  "
$ws.Range("D45").Value = "'0.0932"
$ws.Range("E45").Value = "  -2.59%  "
$ws.Range("D46").Value = "'0.563"
$ws.Range("E46").Value = "  -1.04%  "
$ws.Range("E47").Value = "  -3.96%  "
$ws.Range("D48").Value = "'18.33"
$ws.Range("E48").Value = "  -7.32%  "
$ws.Range("D49").Value = "0.0₆0217"
$ws.Range("E49").Value = "  +17.34%  "
$ws.Range("E50").Value = "  -2.46%  "
$ws.Range("E51").Value = "  -0.53%  "
